# A new weekly price record (date 2022-03-03, serial 44623) was inserted as
# the second data row (row 5) of the sheet. All the previously existing
# records at rows 5..23 shift down by one (to rows 6..24); the sheet's used
# range grows from A1:T23 to A1:T24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 5 (and everything below it) down one row, opening up a
# blank row 5 for the new record.
$ws.Rows("5:5").Insert()

# New record's values, in column order A..T.
$newRow = @(
    6,
    "Mercado Mayorista Lo Valledor de Santiago",
    "Metropolitana",
    44623,
    13,
    "Fruta",
    100101,
    "Berries",
    100101006,
    "Higo",
    "Sin especificar",
    "Segunda",
    30,
    16000,
    16000,
    16000,
    '$/bandeja 7 kilos',
    "Región Metropolitana",
    2286,
    7
)

for ($i = 0; $i -lt $newRow.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $newRow[$i]
}
